$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("B4").Value = 2099376
$ws.Range("C4").Value = 9675
$ws.Range("D4").Value = 818235
$ws.Range("E4").Value = 1164872
$ws.Range("G4").Value = 235
$ws.Range("H4").Value = 116269

# Row 5
$ws.Range("B5").Value = 809398
$ws.Range("C5").Value = 3749
$ws.Range("E5").Value = 371544
$ws.Range("G5").Value = 104
$ws.Range("H5").Value = 41162

# Row 7
$ws.Range("B7").Value = 307039
$ws.Range("C7").Value = 8756
$ws.Range("D7").Value = 153224
$ws.Range("E7").Value = 145063
$ws.Range("G7").Value = 251
$ws.Range("H7").Value = 8752

# Row 9
$ws.Range("B9").Value = 290289
$ws.Range("C9").Value = 502

# Row 10
$ws.Range("B10").Value = 236305
$ws.Range("C10").Value = 163
$ws.Range("D10").Value = 173085
$ws.Range("E10").Value = 28997
$ws.Range("G10").Value = 56
$ws.Range("H10").Value = 34223

# Row 12
$ws.Range("B12").Value = 186933
$ws.Range("C12").Value = 138
$ws.Range("E12").Value = 6480

# Row 15
$ws.Range("A15").Value = "Chile"
$ws.Range("B15").Value = 160846
$ws.Range("C15").Value = 6754
$ws.Range("D15").Value = 131358
$ws.Range("E15").Value = 26618
$ws.Range("G15").Value = 222
$ws.Range("H15").Value = 2870

# Row 16
$ws.Range("A16").Value = "Francia"
$ws.Range("B16").Value = 155561
$ws.Range("D16").Value = 72149
$ws.Range("E16").Value = 54066
$ws.Range("H16").Value = 29346

# Row 20
$ws.Range("B20").Value = 97893
$ws.Range("C20").Value = 363
$ws.Range("D20").Value = 58484
$ws.Range("E20").Value = 31361
$ws.Range("G20").Value = 54
$ws.Range("H20").Value = 8048

# Row 32
$ws.Range("D32").Value = 28040
$ws.Range("E32").Value = 11785

# Row 40
$ws.Range("D40").Value = 8743
$ws.Range("E40").Value = 17858
$ws.Range("G40").Value = 7
$ws.Range("H40").Value = 772

# Row 41
$ws.Range("B41").Value = 25250
$ws.Range("C41").Value = 12
$ws.Range("E41").Value = 847
$ws.Range("G41").Value = 2
$ws.Range("H41").Value = 1705

# Row 44
$ws.Range("B44").Value = 22008
$ws.Range("C44").Value = 571
$ws.Range("D44").Value = 12754
$ws.Range("E44").Value = 8686
$ws.Range("G44").Value = 7
$ws.Range("H44").Value = 568

# Row 62
$ws.Range("B62").Value = 10698
$ws.Range("C62").Value = 109
$ws.Range("D62").Value = 7322
$ws.Range("E62").Value = 2625
$ws.Range("G62").Value = 10
$ws.Range("H62").Value = 751

# Row 82
$ws.Range("A82").Value = "Luxemburgo"
$ws.Range("B82").Value = 4055
$ws.Range("C82").Value = 3
$ws.Range("D82").Value = 3918
$ws.Range("E82").Value = 27
$ws.Range("G82").Value = 0
$ws.Range("H82").Value = 110

# Row 83
$ws.Range("A83").Value = "Hungria"
$ws.Range("B83").Value = 4053
$ws.Range("C83").Value = 14
$ws.Range("D83").Value = 2447
$ws.Range("E83").Value = 1051
$ws.Range("G83").Value = 2
$ws.Range("H83").Value = 555

# Row 90
$ws.Range("B90").Value = 3108
$ws.Range("C90").Value = 20
$ws.Range("E90").Value = 1551

# Row 125
$ws.Range("A125").Value = "Guayana Francesa"
$ws.Range("B125").Value = 1043
$ws.Range("C125").Value = 126
$ws.Range("D125").Value = 489
$ws.Range("E125").Value = 552
$ws.Range("H125").Value = 2

# Row 126
$ws.Range("A126").Value = "Republica de Chipre"
$ws.Range("B126").Value = 975
$ws.Range("D126").Value = 807
$ws.Range("E126").Value = 150
$ws.Range("H126").Value = 18

# Row 127
$ws.Range("A127").Value = "Niger"
$ws.Range("B127").Value = 974
$ws.Range("D127").Value = 878
$ws.Range("E127").Value = 31
$ws.Range("H127").Value = 65

# Row 128
$ws.Range("A128").Value = "Jordania"
$ws.Range("B128").Value = 915
$ws.Range("C128").Value = 25
$ws.Range("D128").Value = 671
$ws.Range("E128").Value = 235
$ws.Range("H128").Value = 9

# Row 129
$ws.Range("A129").Value = "Burkina Faso"
$ws.Range("B129").Value = 892
$ws.Range("D129").Value = 791
$ws.Range("E129").Value = 48
$ws.Range("H129").Value = 53
